$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '71.007.07'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +6.00%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.679.98'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +18.73%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '620.02'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +8.07%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '182.72'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +3.32%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.678.03'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +18.66%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.541'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +5.89%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.163'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +8.03%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.64'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +4.65%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.502'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +7.44%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '40.46'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +12.34%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000254'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +6.44%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.290.50'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +18.76%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.684.07'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +18.96%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '71.048.92'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +6.29%  '
$ws.Range('E18').Value = '  +1.72%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.54'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +7.56%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '519.36'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +8.50%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '16.89'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.68%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.31'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +19.51%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.743'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +8.03%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.55'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +13.27%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '88.50'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +6.11%  '
$ws.Range('E26').Value = '  +7.84%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '11.17'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +10.69%  '
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.54'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +10.35%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.18'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +3.53%  '
$ws.Range('E31').Value = '  +12.51%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '31.65'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +13.10%  '
$ws.Range('B33').Value = 'PEPE'
$ws.Range('C33').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0000111'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +18.39%  '
$ws.Range('E34').Value = '  +4.29%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '6.13'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +9.86%  '
$ws.Range('E37').Value = '  +9.39%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.349'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +12.17%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.24'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +11.36%  '
$ws.Range('E40').Value = '  +7.18%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '51.53'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +5.03%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '45.58'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -5.80%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '434.32'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +16.93%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.81'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +6.04%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.119.12'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +11.55%  '
$ws.Range('E46').Value = '  +5.12%  '
$ws.Range('E47').Value = '  +7.32%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '28.36'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +11.42%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '140.60'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +3.35%  '
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('E51').Value = '  +8.99%  '
